$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing C2 value
$ws.Range("C2").Value = 2107

# Data rows: A=index, B=index, C=value, starting at row 3
$values = @(
    @(2, 2, 2113),
    @(3, 3, 2117),
    @(4, 4, 2120),
    @(5, 5, 2121),
    @(6, 6, 2128),
    @(7, 7, 2134),
    @(8, 8, 2136),
    @(9, 9, 2141),
    @(10, 10, 2146)
)

# Copy the formatting of A2 down through A11 so the new index cells match
$ws.Range("A2").Copy()
$ws.Range("A3:A11").PasteSpecial(-4122)

$row = 3
foreach ($entry in $values) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
